# Regenerated BoM (KiBot) update for pedalboard-led-ring:
#  - Capacitor group (row 9): footprint placement changed (X/Y/Rot)
#  - LED group (row 10): part changed from SK6812MINI (LED lib) to
#    SK6812 (local pedalboard-led-ring lib), new datasheet URL, and
#    footprint rotation changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BoM")

# --- Row 9: Capacitor (C1..C12) footprint placement ---
# These are stored as text in the sheet (e.g. "10.3500"), not numbers,
# so force text entry with a leading apostrophe to avoid Excel's
# automatic number conversion.
$ws.Range("O9").Value = "'8.9400"
$ws.Range("P9").Value = "'1.3500"
$ws.Range("Q9").Value = "'180.0000"

# --- Row 10: LED (D1..D12) part/library/datasheet/rotation ---
$ws.Range("C10").Value = "SK6812"
$ws.Range("D10").Value = "pedalboard-led-ring"
$ws.Range("L10").Value = "https://cdn-shop.adafruit.com/product-files/1138/SK6812+LED+datasheet+.pdf"
$ws.Range("Q10").Value = "'180.0000"
